$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (GitHub Actions refresh).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.113.45'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.51%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.680.91'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.89%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.508'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.07%  '

$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.253'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.20%  '

$ws.Range("E9").Value = '  +2.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.24%  '

$ws.Range("E11").Value = '  +4.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.918.46'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.91%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.680.90'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.90%  '

$ws.Range("E14").Value = '  +1.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.525'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.06%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.146.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '238.36'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0739'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.33%  '

$ws.Range("B23").Value = 'Avalanche'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.66%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.67'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.92%  '

$ws.Range("E26").Value = '  +1.89%  '

$ws.Range("E27").Value = '  +1.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.92%  '

$ws.Range("E29").Value = '  -0.20%  '

$ws.Range("E30").Value = '  +0.48%  '

$ws.Range("E31").Value = '  +2.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.487.93'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.11'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.82%  '

$ws.Range("E35").Value = '  +5.39%  '

$ws.Range("E36").Value = '  -0.38%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.909'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.39%  '

$ws.Range("E38").Value = '  +1.77%  '

$ws.Range("E39").Value = '  +2.62%  '

$ws.Range("E40").Value = '  +2.69%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '67.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.28%  '

$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.28'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.71%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.983'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.827.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.779'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.23%  '

$ws.Range("E47").Value = '  +0.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.54'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.87%  '

$ws.Range("E49").Value = '  +5.04%  '

$ws.Range("E50").Value = '  +1.10%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.24%  '
